$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4: column B/C/D now reference shifted shared strings
# (a new string "Inflammatory-Mac" was inserted before "Rbp4"), and several
# numeric values were recalculated using new TPM-based figures.

# Row 2
$ws.Range("B2").Value = "Rbp4"
$ws.Range("C2").Value = "Stra6"
$ws.Range("D2").Value = "ECs"
$ws.Range("I2").Value = 0.93806027808652
$ws.Range("J2").Value = 0.93806027808652
$ws.Range("O2").Value = 0.01654931057352943
$ws.Range("P2").Value = 0.01654931057352943
$ws.Range("S2").Value = 0.0155242508787452
$ws.Range("T2").Value = 0.0155242508787452

# Row 3
$ws.Range("B3").Value = "Rbp4"
$ws.Range("C3").Value = "Stra6"
$ws.Range("I3").Value = 0.93806027808652
$ws.Range("J3").Value = 0.93806027808652
$ws.Range("O3").Value = 0.3730527584747022
$ws.Range("P3").Value = 0.3730527584747023
$ws.Range("S3").Value = 0.3499459743557226
$ws.Range("T3").Value = 0.3499459743557226

# Row 4
$ws.Range("B4").Value = "Rbp4"
$ws.Range("C4").Value = "Stra6"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 0.93806027808652
$ws.Range("J4").Value = 0.93806027808652
$ws.Range("M4").Value = 0.9924266666666667
$ws.Range("N4").Value = 2.97728
$ws.Range("O4").Value = 0.6103979309517683
$ws.Range("P4").Value = 0.6103979309517683
$ws.Range("Q4").Value = 0.6960357961955554
$ws.Range("R4").Value = 6.264322165759999
$ws.Range("S4").Value = 0.5725900528520522
$ws.Range("T4").Value = 0.5725900528520522

# New row 5 - Inflammatory-Mac -> Rbp4/Stra6 -> ECs
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Rbp4"
$ws.Range("C5").Value = "Stra6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04630966666666667
$ws.Range("H5").Value = 0.138929
$ws.Range("I5").Value = 0.06193972191347994
$ws.Range("J5").Value = 0.06193972191347993
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.026907
$ws.Range("N5").Value = 0.080721
$ws.Range("O5").Value = 0.01654931057352943
$ws.Range("P5").Value = 0.01654931057352943
$ws.Range("Q5").Value = 0.001246054201
$ws.Range("R5").Value = 0.011214487809
$ws.Range("S5").Value = 0.001025059694784226
$ws.Range("T5").Value = 0.001025059694784226

# New row 6 - Inflammatory-Mac -> Rbp4/Stra6 -> FAPs
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Rbp4"
$ws.Range("C6").Value = "Stra6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04630966666666667
$ws.Range("H6").Value = 0.138929
$ws.Range("I6").Value = 0.06193972191347994
$ws.Range("J6").Value = 0.06193972191347993
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6065346666666667
$ws.Range("N6").Value = 1.819604
$ws.Range("O6").Value = 0.3730527584747022
$ws.Range("P6").Value = 0.3730527584747023
$ws.Range("Q6").Value = 0.02808841823511111
$ws.Range("R6").Value = 0.252795764116
$ws.Range("S6").Value = 0.02310678411897965
$ws.Range("T6").Value = 0.02310678411897965

# New row 7 - Inflammatory-Mac -> Rbp4/Stra6 -> MuSCs
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Rbp4"
$ws.Range("C7").Value = "Stra6"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04630966666666667
$ws.Range("H7").Value = 0.138929
$ws.Range("I7").Value = 0.06193972191347994
$ws.Range("J7").Value = 0.06193972191347993
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9924266666666667
$ws.Range("N7").Value = 2.97728
$ws.Range("O7").Value = 0.6103979309517683
$ws.Range("P7").Value = 0.6103979309517683
$ws.Range("Q7").Value = 0.04595894812444445
$ws.Range("R7").Value = 0.41363053312
$ws.Range("S7").Value = 0.03780787809971606
$ws.Range("T7").Value = 0.03780787809971606
